$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "'GJ2010"
$ws.Range("B5").Value = "'ARRIVAL"
$ws.Range("C5").Value = "'Potassium Chloride"
$ws.Range("D5").Value = "'0.12"
$ws.Range("E5").Value = "'2026-02-17"
$ws.Range("F5").Value = "'17:30"
$ws.Range("G5").Value = "'B101003"
$ws.Range("H5").Value = "'O101"
$ws.Range("I5").Value = "'abc"
$ws.Range("J5").Value = "'Central Admin"
$ws.Range("K5").Value = "'16-02-2026 17:27"

$ws.Range("A5:K5").Style = "Normal"
